$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: new test case "tc007" (create interest group) ---
# Clone formatting from row 7 (A = s1 style, B:D = s2 style) before writing values,
# so the new row matches the look of the existing data rows.
$ws.Range("A7:D7").Copy()
$ws.Range("A8:D8").PasteSpecial(-4122)

$ws.Range("A8").Value = "tc007"
$ws.Range("B8").Value = "创建兴趣小组成功"
$ws.Range("C8").Value = "create_interest_group_success_p"
$ws.Range("D8").Value = "create_interest_group_success_e"

# --- Rows 9-30: placeholder cases tc008..tc029, case_id filled, rest blank ---
$ws.Range("A7").Copy()
$ws.Range("A9:A30").PasteSpecial(-4122)

$ws.Range("B4:D4").Copy()
$ws.Range("B9:D30").PasteSpecial(-4122)

For ($r = 9; $r -le 30; $r++) {
    $n = $r - 1
    $ws.Range("A$r").Value = "tc{0:D3}" -f $n
}

# --- Column widths: widen C and D to fit the longer case names ---
$ws.Columns("C").ColumnWidth = 28.8
$ws.Columns("D").ColumnWidth = 28.8

# --- Selection moves to B24 ---
$ws.Range("B24").Select() | Out-Null

$excel.CutCopyMode = $false
